# Insert a new daily price record for "Vega Monumental Concepción - Ají"
# at row 36, pushing the existing rows 36:81 down to 37:82 (weekly refresh
# of the consolidated series).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data (and formatting) down by one row, starting at row 36.
$ws.Rows("36:36").Insert()

# Populate the newly inserted row with the latest observation.
$ws.Cells.Item(36, 1).Value = 11
$ws.Cells.Item(36, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(36, 3).Value = "Bíobío"
$ws.Cells.Item(36, 4).Value = 44629
$ws.Cells.Item(36, 5).Value = 8
$ws.Cells.Item(36, 6).Value = 100112021
$ws.Cells.Item(36, 7).Value = "Ají"
$ws.Cells.Item(36, 8).Value = "Americana (o)"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 35
$ws.Cells.Item(36, 11).Value = 20000
$ws.Cells.Item(36, 12).Value = 21000
$ws.Cells.Item(36, 13).Value = 20571
$ws.Cells.Item(36, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(36, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(36, 16).Value = 823
$ws.Cells.Item(36, 17).Value = 25
$ws.Cells.Item(36, 18).Value = "Hortaliza"
